$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data (rows 2-9), replacing/extending the original 3 data rows (rows 2-4)
$data = @(
    @{ Row = 2; Text = "biophysicist "; X = 463; Y = 669;    W = 105; H = 23 },
    @{ Row = 3; Text = "double ";       X = 564; Y = 694;    W = 63;  H = 23 },
    @{ Row = 4; Text = "helix.";        X = 627; Y = 694;    W = 45;  H = 23 },
    @{ Row = 5; Text = "3D ";           X = 168; Y = 851.8;  W = 30;  H = 23 },
    @{ Row = 6; Text = "model ";        X = 198; Y = 851.8;  W = 58;  H = 23 },
    @{ Row = 7; Text = "human ";        X = 350; Y = 1176;   W = 64;  H = 23 },
    @{ Row = 8; Text = "genome, ";      X = 414; Y = 1176;   W = 80;  H = 23 },
    @{ Row = 9; Text = "discrimination. "; X = 243; Y = 1242.4; W = 128; H = 23 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Text
    $ws.Cells.Item($r, 2).Value = $item.X
    $ws.Cells.Item($r, 3).Value = $item.Y
    $ws.Cells.Item($r, 4).Value = $item.W
    $ws.Cells.Item($r, 5).Value = $item.H
}
